$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts "Data" and every later column
# one position to the right)
$ws.Range("D1").EntireColumn.Insert()

# New header: "Ano" takes the freed-up column D
$ws.Range("D1").Value = "Ano"

# Match the column width Excel would inherit from column C on insert
$ws.Range("D1").ColumnWidth = 10.5

# Leave the active cell on the new header, like the authored workbook
$ws.Range("D1").Select() | Out-Null
